# Dodati su testovii za otvaranje svih stranica
# Add the "Terms" list of values to the "Other" sheet and make it the active tab.

$wb = $excel.ActiveWorkbook

# "Other" is the third worksheet (URL, Users, Other)
$wsOther = $wb.Worksheets.Item(3)

# Fill A1:A5 with the new terms, backed by new shared strings.
$wsOther.Range("A1").Value = "Terms"
$wsOther.Range("A2").Value = "Git"
$wsOther.Range("A3").Value = "Java"
$wsOther.Range("A4").Value = "JS"
$wsOther.Range("A5").Value = "Web"

# Make "Other" the active/selected sheet and select cell K6 on it,
# which also updates workbookView's activeTab and removes tabSelected
# from the previously active "URL" sheet.
$wsOther.Activate()
$wsOther.Range("K6").Select()
